$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("random_forest")

# New configuration labels for column A (rows 2-81), restoring the
# original shared-string text so Excel dedupes against the strings
# already used on the other sheets instead of keeping 40 stale duplicates.
$configValues = @(
    'config=4,  shuffle=False, pca=False, scal=True, minmax=False, lasso=False',
    'config=4,  shuffle=False, pca=False, scal=False, minmax=True, lasso=False',
    'config=4,  shuffle=False, pca=False, scal=False, minmax=False, lasso=False',
    'config=4,  shuffle=True, pca=False, scal=False, minmax=True, lasso=True',
    'config=1,  shuffle=False, pca=False, scal=True, minmax=False, lasso=True',
    'config=1,  shuffle=False, pca=False, scal=False, minmax=True, lasso=True',
    'config=2,  shuffle=False, pca=False, scal=True, minmax=False, lasso=True',
    'config=2,  shuffle=False, pca=False, scal=False, minmax=True, lasso=True',
    'config=4,  shuffle=False, pca=True, scal=False, minmax=True, lasso=False',
    'config=4,  shuffle=True, pca=True, scal=True, minmax=False, lasso=False',
    'config=1,  shuffle=True, pca=False, scal=True, minmax=False, lasso=True',
    'config=2,  shuffle=True, pca=False, scal=True, minmax=False, lasso=True',
    'config=1,  shuffle=False, pca=True, scal=True, minmax=False, lasso=True',
    'config=2,  shuffle=False, pca=True, scal=True, minmax=False, lasso=True',
    'config=4,  shuffle=True, pca=True, scal=False, minmax=True, lasso=False',
    'config=4,  shuffle=False, pca=True, scal=True, minmax=False, lasso=False',
    'config=4,  shuffle=False, pca=False, scal=True, minmax=False, lasso=True',
    'config=1,  shuffle=False, pca=True, scal=True, minmax=False, lasso=False',
    'config=1,  shuffle=False, pca=True, scal=False, minmax=True, lasso=True',
    'config=1,  shuffle=False, pca=True, scal=False, minmax=True, lasso=False',
    'config=1,  shuffle=False, pca=True, scal=False, minmax=False, lasso=False',
    'config=1,  shuffle=True, pca=False, scal=False, minmax=True, lasso=True',
    'config=2,  shuffle=False, pca=True, scal=True, minmax=False, lasso=False',
    'config=2,  shuffle=False, pca=True, scal=False, minmax=True, lasso=True',
    'config=2,  shuffle=False, pca=True, scal=False, minmax=True, lasso=False',
    'config=2,  shuffle=False, pca=True, scal=False, minmax=False, lasso=False',
    'config=2,  shuffle=False, pca=False, scal=True, minmax=False, lasso=False',
    'config=2,  shuffle=True, pca=False, scal=False, minmax=True, lasso=True',
    'config=2,  shuffle=False, pca=False, scal=False, minmax=True, lasso=False',
    'config=2,  shuffle=False, pca=False, scal=False, minmax=False, lasso=False',
    'config=3,  shuffle=False, pca=True, scal=True, minmax=False, lasso=True',
    'config=3,  shuffle=False, pca=True, scal=True, minmax=False, lasso=False',
    'config=3,  shuffle=False, pca=True, scal=False, minmax=True, lasso=True',
    'config=3,  shuffle=False, pca=True, scal=False, minmax=True, lasso=False',
    'config=3,  shuffle=True, pca=True, scal=False, minmax=False, lasso=False',
    'config=3,  shuffle=False, pca=True, scal=False, minmax=False, lasso=False',
    'config=3,  shuffle=True, pca=False, scal=True, minmax=False, lasso=False',
    'config=3,  shuffle=True, pca=False, scal=False, minmax=True, lasso=False',
    'config=3,  shuffle=True, pca=False, scal=False, minmax=False, lasso=False',
    'config=4,  shuffle=True, pca=False, scal=True, minmax=False, lasso=False',
    'config=4,  shuffle=True, pca=False, scal=False, minmax=True, lasso=False',
    'config=4,  shuffle=True, pca=False, scal=False, minmax=False, lasso=False',
    'config=4,  shuffle=True, pca=False, scal=True, minmax=False, lasso=True',
    'config=4,  shuffle=True, pca=True, scal=False, minmax=True, lasso=True',
    'config=2,  shuffle=True, pca=True, scal=False, minmax=False, lasso=False',
    'config=3,  shuffle=False, pca=False, scal=True, minmax=False, lasso=False',
    'config=3,  shuffle=False, pca=False, scal=False, minmax=True, lasso=True',
    'config=3,  shuffle=False, pca=False, scal=False, minmax=True, lasso=False',
    'config=3,  shuffle=False, pca=False, scal=False, minmax=False, lasso=False',
    'config=1,  shuffle=True, pca=True, scal=True, minmax=False, lasso=False',
    'config=1,  shuffle=True, pca=True, scal=False, minmax=False, lasso=False',
    'config=1,  shuffle=True, pca=False, scal=True, minmax=False, lasso=False',
    'config=1,  shuffle=True, pca=False, scal=False, minmax=True, lasso=False',
    'config=1,  shuffle=True, pca=False, scal=False, minmax=False, lasso=False',
    'config=2,  shuffle=True, pca=True, scal=True, minmax=False, lasso=False',
    'config=2,  shuffle=True, pca=False, scal=True, minmax=False, lasso=False',
    'config=2,  shuffle=True, pca=False, scal=False, minmax=True, lasso=False',
    'config=2,  shuffle=True, pca=False, scal=False, minmax=False, lasso=False',
    'config=3,  shuffle=True, pca=True, scal=False, minmax=True, lasso=False',
    'config=3,  shuffle=True, pca=False, scal=False, minmax=True, lasso=True',
    'config=4,  shuffle=False, pca=True, scal=True, minmax=False, lasso=True',
    'config=4,  shuffle=True, pca=True, scal=False, minmax=False, lasso=False',
    'config=4,  shuffle=False, pca=True, scal=False, minmax=True, lasso=True',
    'config=1,  shuffle=False, pca=False, scal=True, minmax=False, lasso=False',
    'config=1,  shuffle=False, pca=False, scal=False, minmax=True, lasso=False',
    'config=1,  shuffle=False, pca=False, scal=False, minmax=False, lasso=False',
    'config=3,  shuffle=False, pca=False, scal=True, minmax=False, lasso=True',
    'config=4,  shuffle=False, pca=True, scal=False, minmax=False, lasso=False',
    'config=1,  shuffle=True, pca=True, scal=False, minmax=True, lasso=True',
    'config=1,  shuffle=True, pca=True, scal=False, minmax=True, lasso=False',
    'config=2,  shuffle=True, pca=True, scal=False, minmax=True, lasso=True',
    'config=2,  shuffle=True, pca=True, scal=False, minmax=True, lasso=False',
    'config=3,  shuffle=True, pca=True, scal=True, minmax=False, lasso=False',
    'config=3,  shuffle=True, pca=False, scal=True, minmax=False, lasso=True',
    'config=4,  shuffle=False, pca=False, scal=False, minmax=True, lasso=True',
    'config=1,  shuffle=True, pca=True, scal=True, minmax=False, lasso=True',
    'config=2,  shuffle=True, pca=True, scal=True, minmax=False, lasso=True',
    'config=4,  shuffle=True, pca=True, scal=True, minmax=False, lasso=True',
    'config=3,  shuffle=True, pca=True, scal=True, minmax=False, lasso=True',
    'config=3,  shuffle=True, pca=True, scal=False, minmax=True, lasso=True'
)

for ($i = 0; $i -lt $configValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $configValues[$i]
}

# Selection moved from F16 to A9
$ws.Range("A9").Select()

# Column A widened from 53.66 to (as close as the host allows to) 63.21875
$ws.Columns.Item(1).ColumnWidth = 62.333333333333336

